# Add files via upload
# Append the next batch of ticker codes below the existing list and bring
# the sheet's row heights / column widths in line with the refreshed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting already used by the existing entries (A2:A4) down onto
# the new rows before filling in their values.
$ws.Range("A2:A4").Copy()
$ws.Range("A5:A17").PasteSpecial(-4122)

$values = @(5406, 6330, 1662, 7202, 4004, 6315, 5016, 5892, 4042, 5713, 5214, 7606, 6902)
$r = 5
foreach ($v in $values) {
    $ws.Cells.Item($r, 1).Value = $v
    $r = $r + 1
}

# The refreshed sheet carries an explicit 15.75pt row height on every row
# down through 100 (not just the rows holding data).
for ($i = 1; $i -le 100; $i++) {
    $ws.Rows.Item($i).RowHeight = 15.75
}

# Columns A:K get an explicit width as well.
$ws.Range("A1:K1").EntireColumn.ColumnWidth = 11.7

Write-Output "applied ticker refresh"
